$d = $word.ActiveDocument

$d.Content.Find.Execute(
  "ResourceMonad<ResourceClass : SubjectResource, etc.>",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "ResourceMonad<ResourceClass : IResource / Resource URNs, etc.>", 2)

$d.Content.Find.Execute(
  "KindMonad<KindClass, SubjectKind, etc.> Monad",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "KindMonad<KindClass :  SubjectKind, etc.> Monad", 2)

$d.Content.Find.Execute(
  "ResourceQuadMonad<ResourceQuadClass : ISubject, etc.> Monad",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "ResourceQuadMonad<ResourceQuadClass : SubjectResource, etc.> Monad", 2)

$d.Content.Find.Execute(
  "OccurrenceQuadMonad<OccurrenceQuadClass : ISubject, etc.> Monad",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "OccurrenceQuadMonad<OccurrenceQuadClass : Subject, etc.> Monad", 2)

$d.Content.Find.Execute(
  "IResource : IQuad",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "IResource : (ResourceURN Resources, Occurrences, Kinds Bindings) : IQuad", 2)
